# Update the Lafayette_B team-specific transition-probability matrix with
# refreshed figures after adding more simulated games (recomputed ratios).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (state 0)
$ws.Range("B2").Value = 0.2035087719298246
$ws.Range("C2").Value = 0.5403508771929825
$ws.Range("J2").Value = 0.01052631578947368
$ws.Range("P2").Value = 0.143859649122807
$ws.Range("S2").Value = 0.1017543859649123

# Row 3 (state 1)
$ws.Range("C3").Value = 0.04242424242424243
$ws.Range("P3").Value = 0.7878787878787878
$ws.Range("S3").Value = 0.1696969696969697

# Row 4 (state 2)
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.3333333333333333

# Row 6 (state 4)
$ws.Range("B6").Value = 0.05286343612334802
$ws.Range("E6").Value = 0.004405286343612335
$ws.Range("F6").Value = 0.08370044052863436
$ws.Range("J6").Value = 0.2378854625550661
$ws.Range("O6").Value = 0.02643171806167401
$ws.Range("Q6").Value = 0.1409691629955947
$ws.Range("R6").Value = 0.04405286343612335
$ws.Range("S6").Value = 0.4096916299559472

# Row 7 (state 5)
$ws.Range("B7").Value = 0.124031007751938
$ws.Range("D7").Value = 0.02325581395348837
$ws.Range("F7").Value = 0.04651162790697674
$ws.Range("J7").Value = 0.124031007751938
$ws.Range("O7").Value = 0.007751937984496124
$ws.Range("Q7").Value = 0.1472868217054264
$ws.Range("R7").Value = 0.08527131782945736
$ws.Range("S7").Value = 0.4418604651162791

# Row 8 (state 6)
$ws.Range("B8").Value = 0.08780487804878048
$ws.Range("D8").Value = 0.01219512195121951
$ws.Range("E8").Value = 0.002439024390243902
$ws.Range("F8").Value = 0.04878048780487805
$ws.Range("J8").Value = 0.1365853658536585
$ws.Range("O8").Value = 0.03170731707317073
$ws.Range("Q8").Value = 0.1731707317073171
$ws.Range("R8").Value = 0.1048780487804878
$ws.Range("S8").Value = 0.4024390243902439

# Row 9 (state 7)
$ws.Range("B9").Value = 0.08762886597938144
$ws.Range("D9").Value = 0.04639175257731959
$ws.Range("E9").Value = 0.005154639175257732
$ws.Range("F9").Value = 0.08247422680412371
$ws.Range("J9").Value = 0.09278350515463918
$ws.Range("O9").Value = 0.02577319587628866
$ws.Range("Q9").Value = 0.134020618556701
$ws.Range("R9").Value = 0.134020618556701
$ws.Range("S9").Value = 0.3917525773195876

# Row 10 (state 8)
$ws.Range("B10").Value = 0.1076335877862595
$ws.Range("D10").Value = 0.02061068702290076
$ws.Range("E10").Value = 0.0007633587786259542
$ws.Range("F10").Value = 0.07862595419847328
$ws.Range("J10").Value = 0.1145038167938931
$ws.Range("O10").Value = 0.02366412213740458
$ws.Range("Q10").Value = 0.1961832061068702
$ws.Range("R10").Value = 0.08091603053435115
$ws.Range("S10").Value = 0.3770992366412214

# Row 11 (state 9)
$ws.Range("G11").Value = 0.1357466063348416
$ws.Range("J11").Value = 0.09049773755656108
$ws.Range("K11").Value = 0.2036199095022624
$ws.Range("L11").Value = 0.5520361990950227
$ws.Range("S11").Value = 0.01809954751131222

# Row 12 (state 10)
$ws.Range("G12").Value = 0.6885245901639344
$ws.Range("J12").Value = 0.2377049180327869
$ws.Range("K12").Value = 0.00819672131147541
$ws.Range("L12").Value = 0.02459016393442623
$ws.Range("S12").Value = 0.04098360655737705

# Row 13 (state 11)
$ws.Range("G13").Value = 0.625
$ws.Range("J13").Value = 0.34375
$ws.Range("S13").Value = 0.03125

# Row 15 (state 13)
$ws.Range("F15").Value = 0.007905138339920948
$ws.Range("H15").Value = 0.1462450592885375
$ws.Range("I15").Value = 0.05928853754940711
$ws.Range("J15").Value = 0.4229249011857708
$ws.Range("K15").Value = 0.04347826086956522
$ws.Range("M15").Value = 0.003952569169960474
$ws.Range("O15").Value = 0.06719367588932806
$ws.Range("S15").Value = 0.2490118577075099

# Row 16 (state 14)
$ws.Range("F16").Value = 0.0101010101010101
$ws.Range("H16").Value = 0.1919191919191919
$ws.Range("I16").Value = 0.06060606060606061
$ws.Range("J16").Value = 0.4343434343434344
$ws.Range("K16").Value = 0.09595959595959595
$ws.Range("M16").Value = 0.01515151515151515
$ws.Range("N16").Value = 0.005050505050505051
$ws.Range("O16").Value = 0.05555555555555555
$ws.Range("S16").Value = 0.1313131313131313

# Row 17 (state 15)
$ws.Range("F17").Value = 0.01496259351620948
$ws.Range("H17").Value = 0.1945137157107232
$ws.Range("I17").Value = 0.09725685785536159
$ws.Range("J17").Value = 0.3990024937655861
$ws.Range("K17").Value = 0.08478802992518704
$ws.Range("M17").Value = 0.01496259351620948
$ws.Range("N17").Value = 0.002493765586034913
$ws.Range("O17").Value = 0.0598503740648379
$ws.Range("S17").Value = 0.1321695760598504

# Row 18 (state 16)
$ws.Range("F18").Value = 0.02564102564102564
$ws.Range("H18").Value = 0.1641025641025641
$ws.Range("I18").Value = 0.09743589743589744
$ws.Range("J18").Value = 0.4974358974358974
$ws.Range("K18").Value = 0.04615384615384616
$ws.Range("M18").Value = 0.01538461538461539
$ws.Range("O18").Value = 0.06153846153846154
$ws.Range("S18").Value = 0.09230769230769231

# Row 19 (state 17)
$ws.Range("F19").Value = 0.01883701883701884
$ws.Range("H19").Value = 0.1850941850941851
$ws.Range("I19").Value = 0.09009009009009009
$ws.Range("J19").Value = 0.4266994266994267
$ws.Range("K19").Value = 0.0819000819000819
$ws.Range("M19").Value = 0.0171990171990172
$ws.Range("O19").Value = 0.08517608517608517
$ws.Range("S19").Value = 0.095004095004095
